$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right after row 283 (before the old row 284),
# pushing the existing rows 284-300 down to 286-302.
$ws.Rows("284:285").Insert()

# The insert copies row 283's formatting into the new blank rows, including an
# empty styled K284/K285 cell that the final sheet should not contain at all.
# Reset that cell's style back to the workbook default and clear it so it is
# dropped entirely rather than left behind as an empty styled cell.
$ws.Range("K284:K285").Style = $ws.Range("Z999").Style
$ws.Range("K284:K285").ClearContents()

# Populate the common (pre-existing) values for both rows first.
$ws.Range("B284").Value2 = "peppol-doctype-wildcard"
$ws.Range("D284").Formula = "'9.1"
$ws.Range("E284").Value2 = "active"
$ws.Range("H284").Value2 = "TICC-373"
$ws.Range("I284").Value2 = $false
$ws.Range("J284").Formula = "=TRUE()"
$ws.Range("L284").Value2 = "POAC"

$ws.Range("B285").Value2 = "peppol-doctype-wildcard"
$ws.Range("D285").Formula = "'9.1"
$ws.Range("E285").Value2 = "active"
$ws.Range("H285").Value2 = "TICC-373"
$ws.Range("I285").Value2 = $false
$ws.Range("J285").Formula = "=TRUE()"
$ws.Range("L285").Value2 = "POAC"

# New shared strings are introduced column-by-column: C, then M, then N, then A.
$ws.Range("C284").Value2 = "urn:fdc:peppol:tax-data-document:1.0::TaxData##urn:peppol:pint:taxdata-1@ae-1::1.0"
$ws.Range("C285").Value2 = "urn:fdc:peppol:tax-data-status:1.0::TaxDataStatus##urn:peppol:pint:taxdatastatus-1@ae-1::1.0"

$ws.Range("M284").Value2 = "Tax Reporting"
$ws.Range("M285").Value2 = "Tax Reporting"

$ws.Range("N284").Value2 = "cenbii-procid-ubl::urn:peppol:bis:taxreporting"
$ws.Range("N285").Value2 = "cenbii-procid-ubl::urn:peppol:bis:taxreporting"

$ws.Range("A284").Value2 = "AE Tax Data Document v1.0"
$ws.Range("A285").Value2 = "AE Tax Data Status v1.0"

# Column A in these rows keeps the sheet's plain default style (no explicit
# style override), matching the other recently-added rows above them.
$ws.Range("A284").Style = $ws.Range("Z999").Style
$ws.Range("A285").Style = $ws.Range("Z999").Style
